# edit.ps1 — generated script implementing commit "feat: add 2022-Q1 data"
#
# Summary of the change:
#   * Insert a new worksheet named '2022-Q1' immediately before the existing
#     grand-total sheet '总计', containing the per-fund holding breakdown for
#     2022-Q1 (same column layout as the most recent quarter, '2021-Q4').
#   * Prepend a '2022-Q1' row (34 funds held, 33.92 billion yuan total market
#     value) to the '总计' summary sheet.
#
# Implementation note on sheet identity: to reproduce the target workbook's
# rId/sheetId layout exactly, we repurpose the CURRENT '总计' sheet object
# (which already owns sheetId=6/rId6) as the new '2022-Q1' sheet, and create
# a brand-new sheet to hold the (updated) '总计' content — that new sheet
# naturally receives the next sheetId/rId (7), matching a '总计' sheet that
# was re-created after the insert.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Helpers
# ---------------------------------------------------------------------------

# Write $val into the cell as literal TEXT, never auto-coerced to a number
# (needed for fund codes with leading zeros, and for numeric-looking metrics
# that the source data stores as text). NumberFormat is forced to "@" for
# the assignment so the smart-input parser leaves the value as a string, then
# immediately cleared again so the cell keeps the sheet's default styling.
function Set-TextCell($ws, $row, $col, $val) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.ClearFormats()
}

# Write $val into the cell as a real number.
function Set-NumCell($ws, $row, $col, $val) {
    $ws.Cells.Item($row, $col).Value = $val
}

# ---------------------------------------------------------------------------
# Step 1 — repurpose the existing '总计' sheet object as the new '2022-Q1' tab
# ---------------------------------------------------------------------------

$styleSrc = $wb.Worksheets.Item("2021-Q4")
$q1 = $wb.Worksheets.Item("总计")
$q1.Name = "2022-Q1"
$q1.Cells.Clear()

# Header row (B1:H1), bold/centered/bordered like the other quarterly sheets
$styleSrc.Range("B1:H1").Copy()
$q1.Range("B1:H1").PasteSpecial(-4122)
$q1.Cells.Item(1, 2).Value = "基金代码"
$q1.Cells.Item(1, 3).Value = "基金名称"
$q1.Cells.Item(1, 4).Value = "基金规模"
$q1.Cells.Item(1, 5).Value = "股票总仓位"
$q1.Cells.Item(1, 6).Value = "仓位占比"
$q1.Cells.Item(1, 7).Value = "持有市值(亿元)"
$q1.Cells.Item(1, 8).Value = "仓位排名"

# Row index column (A2:A35), same styling as the other quarterly sheets' index column
$styleSrc.Range("A2:A35").Copy()
$q1.Range("A2:A35").PasteSpecial(-4122)

# Per-fund data rows
Set-NumCell $q1 2 1 0
Set-TextCell $q1 2 2 "004997"
Set-TextCell $q1 2 3 "广发高端制造股票A"
Set-TextCell $q1 2 4 "148.04"
Set-TextCell $q1 2 5 "94.19"
Set-TextCell $q1 2 6 "6.97"
Set-TextCell $q1 2 7 "10.3184"
Set-NumCell $q1 2 8 7
Set-NumCell $q1 3 1 1
Set-TextCell $q1 3 2 "011121"
Set-TextCell $q1 3 3 "广发兴诚混合A"
Set-TextCell $q1 3 4 "53.52"
Set-TextCell $q1 3 5 "94.10"
Set-TextCell $q1 3 6 "9.67"
Set-TextCell $q1 3 7 "5.1754"
Set-NumCell $q1 3 8 3
Set-NumCell $q1 4 1 2
Set-TextCell $q1 4 2 "011479"
Set-TextCell $q1 4 3 "广发诚享混合A"
Set-TextCell $q1 4 4 "44.44"
Set-TextCell $q1 4 5 "93.13"
Set-TextCell $q1 4 6 "10.14"
Set-TextCell $q1 4 7 "4.5062"
Set-NumCell $q1 4 8 1
Set-NumCell $q1 5 1 3
Set-TextCell $q1 5 2 "011130"
Set-TextCell $q1 5 3 "广发兴诚混合C"
Set-TextCell $q1 5 4 "28.81"
Set-TextCell $q1 5 5 "94.10"
Set-TextCell $q1 5 6 "9.67"
Set-TextCell $q1 5 7 "2.7859"
Set-NumCell $q1 5 8 3
Set-NumCell $q1 6 1 4
Set-TextCell $q1 6 2 "161834"
Set-TextCell $q1 6 3 "银华鑫锐灵活配置混合（LOF）"
Set-TextCell $q1 6 4 "67.33"
Set-TextCell $q1 6 5 "81.90"
Set-TextCell $q1 6 6 "3.43"
Set-TextCell $q1 6 7 "2.3094"
Set-NumCell $q1 6 8 2
Set-NumCell $q1 7 1 5
Set-TextCell $q1 7 2 "501022"
Set-TextCell $q1 7 3 "银华鑫盛灵活配置混合（LOF）"
Set-TextCell $q1 7 4 "61.98"
Set-TextCell $q1 7 5 "79.75"
Set-TextCell $q1 7 6 "3.35"
Set-TextCell $q1 7 7 "2.0763"
Set-NumCell $q1 7 8 2
Set-NumCell $q1 8 1 6
Set-TextCell $q1 8 2 "519732"
Set-TextCell $q1 8 3 "交银定期支付双息平衡混合"
Set-TextCell $q1 8 4 "40.83"
Set-TextCell $q1 8 5 "67.67"
Set-TextCell $q1 8 6 "3.53"
Set-TextCell $q1 8 7 "1.4413"
Set-NumCell $q1 8 8 3
Set-NumCell $q1 9 1 7
Set-TextCell $q1 9 2 "001256"
Set-TextCell $q1 9 3 "泓德优选成长混合"
Set-TextCell $q1 9 4 "25.31"
Set-TextCell $q1 9 5 "91.77"
Set-TextCell $q1 9 6 "3.56"
Set-TextCell $q1 9 7 "0.9010"
Set-NumCell $q1 9 8 7
Set-NumCell $q1 10 1 8
Set-TextCell $q1 10 2 "004814"
Set-TextCell $q1 10 3 "中欧红利优享灵活配置混合A"
Set-TextCell $q1 10 4 "22.96"
Set-TextCell $q1 10 5 "93.96"
Set-TextCell $q1 10 6 "2.84"
Set-TextCell $q1 10 7 "0.6521"
Set-NumCell $q1 10 8 7
Set-NumCell $q1 11 1 9
Set-TextCell $q1 11 2 "010160"
Set-TextCell $q1 11 3 "广发高端制造股票C"
Set-TextCell $q1 11 4 "8.51"
Set-TextCell $q1 11 5 "94.19"
Set-TextCell $q1 11 6 "6.97"
Set-TextCell $q1 11 7 "0.5931"
Set-NumCell $q1 11 8 7
Set-NumCell $q1 12 1 10
Set-TextCell $q1 12 2 "180010"
Set-TextCell $q1 12 3 "银华优质增长混合"
Set-TextCell $q1 12 4 "22.68"
Set-TextCell $q1 12 5 "77.41"
Set-TextCell $q1 12 6 "2.18"
Set-TextCell $q1 12 7 "0.4944"
Set-NumCell $q1 12 8 9
Set-NumCell $q1 13 1 11
Set-TextCell $q1 13 2 "011480"
Set-TextCell $q1 13 3 "广发诚享混合C"
Set-TextCell $q1 13 4 "4.52"
Set-TextCell $q1 13 5 "93.13"
Set-TextCell $q1 13 6 "10.14"
Set-TextCell $q1 13 7 "0.4583"
Set-NumCell $q1 13 8 1
Set-NumCell $q1 14 1 12
Set-TextCell $q1 14 2 "011405"
Set-TextCell $q1 14 3 "银华稳健增长一年持有期混合"
Set-TextCell $q1 14 4 "18.05"
Set-TextCell $q1 14 5 "72.76"
Set-TextCell $q1 14 6 "2.18"
Set-TextCell $q1 14 7 "0.3935"
Set-NumCell $q1 14 8 9
Set-NumCell $q1 15 1 13
Set-TextCell $q1 15 2 "012370"
Set-TextCell $q1 15 3 "银华鑫利一年持有期混合型证券投资基金"
Set-TextCell $q1 15 4 "11.01"
Set-TextCell $q1 15 5 "80.06"
Set-TextCell $q1 15 6 "3.36"
Set-TextCell $q1 15 7 "0.3699"
Set-NumCell $q1 15 8 2
Set-NumCell $q1 16 1 14
Set-TextCell $q1 16 2 "004815"
Set-TextCell $q1 16 3 "中欧红利优享灵活配置混合C"
Set-TextCell $q1 16 4 "8.45"
Set-TextCell $q1 16 5 "93.96"
Set-TextCell $q1 16 6 "2.84"
Set-TextCell $q1 16 7 "0.2400"
Set-NumCell $q1 16 8 7
Set-NumCell $q1 17 1 15
Set-TextCell $q1 17 2 "159870"
Set-TextCell $q1 17 3 "鹏华中证细分化工产业主题ETF"
Set-TextCell $q1 17 4 "8.45"
Set-TextCell $q1 17 5 "98.37"
Set-TextCell $q1 17 6 "2.39"
Set-TextCell $q1 17 7 "0.2020"
Set-NumCell $q1 17 8 10
Set-NumCell $q1 18 1 16
Set-TextCell $q1 18 2 "009230"
Set-TextCell $q1 18 3 "鹏华安和混合A"
Set-TextCell $q1 18 4 "14.02"
Set-TextCell $q1 18 5 "34.45"
Set-TextCell $q1 18 6 "1.15"
Set-TextCell $q1 18 7 "0.1612"
Set-NumCell $q1 18 8 9
Set-NumCell $q1 19 1 17
Set-TextCell $q1 19 2 "009667"
Set-TextCell $q1 19 3 "鹏华安庆混合A"
Set-TextCell $q1 19 4 "11.22"
Set-TextCell $q1 19 5 "38.92"
Set-TextCell $q1 19 6 "1.27"
Set-TextCell $q1 19 7 "0.1425"
Set-NumCell $q1 19 8 10
Set-NumCell $q1 20 1 18
Set-TextCell $q1 20 2 "013393"
Set-TextCell $q1 20 3 "信达澳银价值精选混合A"
Set-TextCell $q1 20 4 "3.61"
Set-TextCell $q1 20 5 "81.31"
Set-TextCell $q1 20 6 "2.59"
Set-TextCell $q1 20 7 "0.0935"
Set-NumCell $q1 20 8 8
Set-NumCell $q1 21 1 19
Set-TextCell $q1 21 2 "000805"
Set-TextCell $q1 21 3 "中银新经济灵活配置混合"
Set-TextCell $q1 21 4 "2.75"
Set-TextCell $q1 21 5 "71.79"
Set-TextCell $q1 21 6 "2.64"
Set-TextCell $q1 21 7 "0.0726"
Set-NumCell $q1 21 8 6
Set-NumCell $q1 22 1 20
Set-TextCell $q1 22 2 "001120"
Set-TextCell $q1 22 3 "东方睿鑫热点挖掘灵活配置混合A"
Set-TextCell $q1 22 4 "1.74"
Set-TextCell $q1 22 5 "78.43"
Set-TextCell $q1 22 6 "3.89"
Set-TextCell $q1 22 7 "0.0677"
Set-NumCell $q1 22 8 9
Set-NumCell $q1 23 1 21
Set-TextCell $q1 23 2 "673020"
Set-TextCell $q1 23 3 "西部利得成长精选灵活配置混合"
Set-TextCell $q1 23 4 "1.36"
Set-TextCell $q1 23 5 "91.42"
Set-TextCell $q1 23 6 "4.58"
Set-TextCell $q1 23 7 "0.0623"
Set-NumCell $q1 23 8 4
Set-NumCell $q1 24 1 22
Set-TextCell $q1 24 2 "009231"
Set-TextCell $q1 24 3 "鹏华安和混合C"
Set-TextCell $q1 24 4 "5.33"
Set-TextCell $q1 24 5 "34.45"
Set-TextCell $q1 24 6 "1.15"
Set-TextCell $q1 24 7 "0.0613"
Set-NumCell $q1 24 8 9
Set-NumCell $q1 25 1 23
Set-TextCell $q1 25 2 "163810"
Set-TextCell $q1 25 3 "中银价值精选灵活配置混合"
Set-TextCell $q1 25 4 "1.90"
Set-TextCell $q1 25 5 "76.15"
Set-TextCell $q1 25 6 "3.22"
Set-TextCell $q1 25 7 "0.0612"
Set-NumCell $q1 25 8 5
Set-NumCell $q1 26 1 24
Set-TextCell $q1 26 2 "516120"
Set-TextCell $q1 26 3 "富国中证细分化工产业主题ETF"
Set-TextCell $q1 26 4 "2.32"
Set-TextCell $q1 26 5 "99.12"
Set-TextCell $q1 26 6 "2.42"
Set-TextCell $q1 26 7 "0.0561"
Set-NumCell $q1 26 8 10
Set-NumCell $q1 27 1 25
Set-TextCell $q1 27 2 "006302"
Set-TextCell $q1 27 3 "银华行业轮动混合"
Set-TextCell $q1 27 4 "2.08"
Set-TextCell $q1 27 5 "75.64"
Set-TextCell $q1 27 6 "2.30"
Set-TextCell $q1 27 7 "0.0478"
Set-NumCell $q1 27 8 8
Set-NumCell $q1 28 1 26
Set-TextCell $q1 28 2 "001121"
Set-TextCell $q1 28 3 "东方睿鑫热点挖掘灵活配置混合C"
Set-TextCell $q1 28 4 "1.10"
Set-TextCell $q1 28 5 "78.43"
Set-TextCell $q1 28 6 "3.89"
Set-TextCell $q1 28 7 "0.0428"
Set-NumCell $q1 28 8 9
Set-NumCell $q1 29 1 27
Set-TextCell $q1 29 2 "930602"
Set-TextCell $q1 29 3 "国信价值智选混合型集合资产管理计划"
Set-TextCell $q1 29 4 "0.50"
Set-TextCell $q1 29 5 "67.38"
Set-TextCell $q1 29 6 "6.82"
Set-TextCell $q1 29 7 "0.0341"
Set-NumCell $q1 29 8 3
Set-NumCell $q1 30 1 28
Set-TextCell $q1 30 2 "009668"
Set-TextCell $q1 30 3 "鹏华安庆混合C"
Set-TextCell $q1 30 4 "2.36"
Set-TextCell $q1 30 5 "38.92"
Set-TextCell $q1 30 6 "1.27"
Set-TextCell $q1 30 7 "0.0300"
Set-NumCell $q1 30 8 10
Set-NumCell $q1 31 1 29
Set-TextCell $q1 31 2 "004244"
Set-TextCell $q1 31 3 "东方周期优选灵活配置混合"
Set-TextCell $q1 31 4 "0.62"
Set-TextCell $q1 31 5 "84.13"
Set-TextCell $q1 31 6 "4.52"
Set-TextCell $q1 31 7 "0.0280"
Set-NumCell $q1 31 8 7
Set-NumCell $q1 32 1 30
Set-TextCell $q1 32 2 "005519"
Set-TextCell $q1 32 3 "银华混改红利灵活配置混合"
Set-TextCell $q1 32 4 "0.62"
Set-TextCell $q1 32 5 "72.27"
Set-TextCell $q1 32 6 "2.30"
Set-TextCell $q1 32 7 "0.0143"
Set-NumCell $q1 32 8 5
Set-NumCell $q1 33 1 31
Set-TextCell $q1 33 2 "510290"
Set-TextCell $q1 33 3 "南方上证380ETF"
Set-TextCell $q1 33 4 "1.75"
Set-TextCell $q1 33 5 "99.12"
Set-TextCell $q1 33 6 "0.81"
Set-TextCell $q1 33 7 "0.0142"
Set-NumCell $q1 33 8 9
Set-NumCell $q1 34 1 32
Set-TextCell $q1 34 2 "013394"
Set-TextCell $q1 34 3 "信达澳银价值精选混合C"
Set-TextCell $q1 34 4 "0.37"
Set-TextCell $q1 34 5 "81.31"
Set-TextCell $q1 34 6 "2.59"
Set-TextCell $q1 34 7 "0.0096"
Set-NumCell $q1 34 8 8
Set-NumCell $q1 35 1 33
Set-TextCell $q1 35 2 "519117"
Set-TextCell $q1 35 3 "浦银安盛基本面400指数"
Set-TextCell $q1 35 4 "0.24"
Set-TextCell $q1 35 5 "92.63"
Set-TextCell $q1 35 6 "0.57"
Set-TextCell $q1 35 7 "0.0014"
Set-NumCell $q1 35 8 10

# ---------------------------------------------------------------------------
# Step 2 — create a fresh sheet to hold the updated '总计' (grand-total) data
# ---------------------------------------------------------------------------

$total = $wb.Worksheets.Add($null, $q1)
$total.Name = "总计"

# Match the page margins used throughout the rest of the workbook (inches,
# expressed in points since PageSetup.*Margin is points-denominated: 1in=72pt).
$total.PageSetup.LeftMargin = 54
$total.PageSetup.RightMargin = 54
$total.PageSetup.TopMargin = 72
$total.PageSetup.BottomMargin = 72
$total.PageSetup.HeaderMargin = 36
$total.PageSetup.FooterMargin = 36

# Header row (B1:D1)
$styleSrc.Range("B1:D1").Copy()
$total.Range("B1:D1").PasteSpecial(-4122)
$total.Cells.Item(1, 2).Value = "日期"
$total.Cells.Item(1, 3).Value = "持有数量(只)"
$total.Cells.Item(1, 4).Value = "持有市值(亿元)"

# Row index column (A2:A7)
$styleSrc.Range("A2:A7").Copy()
$total.Range("A2:A7").PasteSpecial(-4122)

# Data rows — 2022-Q1 prepended, existing quarters shifted down by one
Set-NumCell $total 2 1 0
Set-TextCell $total 2 2 "2022-Q1"
Set-NumCell $total 2 3 34
Set-NumCell $total 2 4 33.92
Set-NumCell $total 3 1 1
Set-TextCell $total 3 2 "2021-Q4"
Set-NumCell $total 3 3 33
Set-NumCell $total 3 4 45.93
Set-NumCell $total 4 1 2
Set-TextCell $total 4 2 "2021-Q3"
Set-NumCell $total 4 3 74
Set-NumCell $total 4 4 65.99
Set-NumCell $total 5 1 3
Set-TextCell $total 5 2 "2021-Q2"
Set-NumCell $total 5 3 145
Set-NumCell $total 5 4 86.85
Set-NumCell $total 6 1 4
Set-TextCell $total 6 2 "2021-Q1"
Set-NumCell $total 6 3 123
Set-NumCell $total 6 4 64.21
Set-NumCell $total 7 1 5
Set-TextCell $total 7 2 "2020-Q4"
Set-NumCell $total 7 3 77
Set-NumCell $total 7 4 28.3

Write-Output "edit.ps1 complete"
